# Applies the cryptocurrency price/volume updates described in the commit
# "Updated cryptos list on Tue Mar 28 07:23:28 UTC 2023 with GitHub Actions".
# All the touched cells hold plain text in the workbook (prices/percentages are
# formatted strings, not numbers), so for any value that Excel would otherwise
# auto-convert to a number we briefly force the cell to Text format, assign the
# string, then restore the cell style so no stray number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "27.069.45"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "1.730.34"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "310.42"
$ws.Range("E5").Value = "  -5.24%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.05%  "
Set-TextValue "D7" "0.4880"
$ws.Range("E7").Value = "  +6.57%  "
Set-TextValue "D8" "0.3512"
$ws.Range("E8").Value = "  +0.41%  "
Set-TextValue "D9" "43.49"
$ws.Range("E9").Value = "  +3.74%  "
Set-TextValue "D10" "0.07279"
$ws.Range("E10").Value = "  -1.06%  "
Set-TextValue "D11" "1.050"
$ws.Range("E11").Value = "  -3.18%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.07%  "
Set-TextValue "D13" "20.02"
$ws.Range("E13").Value = "  -2.76%  "
Set-TextValue "D14" "5.889"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "1.733.45"
$ws.Range("E15").Value = "  -1.41%  "
Set-TextValue "D16" "6.900"
$ws.Range("E16").Value = "  -3.78%  "
Set-TextValue "D17" "87.35"
$ws.Range("E17").Value = "  -4.71%  "
Set-TextValue "D18" "0.00001040"
$ws.Range("E18").Value = "  -1.25%  "
Set-TextValue "D19" "0.06414"
$ws.Range("E19").Value = "  -0.12%  "
Set-TextValue "D20" "1.001"
$ws.Range("E20").Value = "  +0.02%  "
Set-TextValue "D21" "16.60"
$ws.Range("E21").Value = "  -1.36%  "
Set-TextValue "D22" "5.695"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "27.110.54"
$ws.Range("E23").Value = "  -2.83%  "
Set-TextValue "D24" "10.94"
$ws.Range("E24").Value = "  -2.08%  "
Set-TextValue "D25" "2.079"
$ws.Range("E25").Value = "  -3.65%  "
Set-TextValue "D26" "153.76"
$ws.Range("E26").Value = "  -4.95%  "
Set-TextValue "D27" "20.01"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "1.929.68"
$ws.Range("E28").Value = "  -1.55%  "
Set-TextValue "D29" "2.089"
$ws.Range("E29").Value = "  -3.02%  "
Set-TextValue "D30" "121.55"
$ws.Range("E30").Value = "  -1.49%  "
Set-TextValue "D31" "1.049"
$ws.Range("E31").Value = "  -2.17%  "
Set-TextValue "D32" "0.09324"
$ws.Range("E32").Value = "  +0.81%  "
Set-TextValue "D33" "3.641"
$ws.Range("E33").Value = "  -0.51%  "
Set-TextValue "D34" "5.383"
$ws.Range("E34").Value = "  -2.89%  "
Set-TextValue "D35" "0.05972"
$ws.Range("E35").Value = "  -2.03%  "
Set-TextValue "D36" "0.02190"
$ws.Range("E36").Value = "  -3.42%  "
Set-TextValue "D37" "1.446"
$ws.Range("E37").Value = "  +5.95%  "
Set-TextValue "D38" "10.99"
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.1998"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D40" "4.774"
$ws.Range("E40").Value = "  -2.57%  "
Set-TextValue "D41" "0.6011"
$ws.Range("E41").Value = "  -2.77%  "
Set-TextValue "D42" "0.9998"
$ws.Range("E42").Value = "  +0.17%  "
Set-TextValue "D43" "1.096"
$ws.Range("E43").Value = "  -6.87%  "
Set-TextValue "D44" "7.511"
$ws.Range("E44").Value = "  -3.28%  "
Set-TextValue "D45" "12.89"
$ws.Range("E45").Value = "  -1.30%  "
Set-TextValue "D46" "3.584"
$ws.Range("E46").Value = "  -3.82%  "
Set-TextValue "D47" "0.5665"
$ws.Range("E47").Value = "  -2.12%  "
Set-TextValue "D48" "118.94"
$ws.Range("E48").Value = "  -3.63%  "
Set-TextValue "D49" "1.852"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("E50").Value = "  -1.50%  "
Set-TextValue "D51" "0.06644"
$ws.Range("E51").Value = "  -2.26%  "
